# Generate Report for Handoff
# Update the "f99c9f5c-8352-4437-9bfc-c98dd3fb24be.md" row (row 3 on every sheet)
# from "In Translation" / "ht" to "Ready for handoff" / "mt", with refreshed
# handoff datetimes, on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# ----- Overview sheet -----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = "Ready for handoff"
$ovw.Range("F3").Value = "Ready for handoff"
$ovw.Range("G3").Value = "2016-08-21 14:13:54"
$ovw.Columns.Item(5).ColumnWidth = 16.333333333333332
$ovw.Columns.Item(6).ColumnWidth = 16.333333333333332

# ----- zh-cn sheet -----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E3").Value = "mt"
$zh.Range("H3").Value = "2016-08-21 14:13:50"
$zh.Columns.Item(3).ColumnWidth = 16.333333333333332

# ----- de-de sheet -----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("E3").Value = "mt"
$de.Range("H3").Value = "2016-08-21 14:13:54"
$de.Columns.Item(3).ColumnWidth = 16.333333333333332
